$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from the last existing header cell (AC1) into the new header cells
# so the new headers match the workbook's bold/centered/bordered header style.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

# Set the new header labels for the season-record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate every data row (2-52) with the team's season record.
for ($r = 2; $r -le 52; $r++) {
    $ws.Cells.Item($r, 30).Value = 93
    $ws.Cells.Item($r, 31).Value = 69
    $ws.Cells.Item($r, 32).Value = 0
}
